$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.060852289199829
$ws.Range("B1").Value = 2.906769275665283
$ws.Range("C1").Value = 5.215261459350586
$ws.Range("D1").Value = 3.573774337768555
$ws.Range("E1").Value = 1.409287571907043
